# Calendrier Sprint 3 - correction de la date de changement automatique.
# Au lieu de calculer la semaine a partir d'aujourd'hui (TODAY()), la
# semaine est desormais figee sur les dates reelles indiquees au debut de
# la Reunion Scrum Master (19/02/2024 - 25/02/2024), et les libelles de
# "Taches terminees le ..." sont decales en consequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Ligne 2 : la date de depart n'est plus TODAY()-1 mais une date fixe
# (lundi 19/02/2024), le reste de la semaine s'en deduit par +1 jour.
$ws.Range("B2").Formula = "=DATE(2024,2,19)"
$ws.Range("C2").Formula = "=B2+1"
$ws.Range("D2:H2").Formula = "=C2+1"

# Libelles "Taches terminees le ..." decales de 4 jours plus tot.
$ws.Range("B19").Value = "Tâches terminé le 20/02/2024"
$ws.Range("B20").Value = "Tâches terminé le 21/02/2024"
$ws.Range("B21").Value = "Tâches terminé le 24/02/2024"
$ws.Range("B22").Value = "Tâches terminé le 25/02/2024"

# Cellule selectionnee au moment de l'enregistrement.
$ws.Range("E12").Select() | Out-Null

$wb.Save()
